$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Reset the lingering selection (M12) left over from the previous save.
$ws.Range("A1").Select()

# Included Multi User test data: two new rows exercising a combined
# "AutoTestAdmin@@AutoTestUser" recipient on the To column, mirroring the
# existing "New Transmittal from Automation" / "Issued for Approval" rows,
# one resolving to Approved and the other to Rejected.

# Row 7: Multi user "To" field, Approved
$ws.Range("A7").Value = "AutoTestAdmin@@AutoTestUser"
$ws.Range("C7").Value = "New Transmittal from Automation"
$ws.Range("D7").Value = "UnTick"
$ws.Range("E7").Value = "Change Note"
$ws.Range("F7").Value = "Issued for Approval"
$ws.Range("L7").Value = "Message for New transmittal"
$ws.Range("M7").Value = "Approved"

# Row 8: Multi user "To" field, Rejected
$ws.Range("A8").Value = "AutoTestAdmin@@AutoTestUser"
$ws.Range("C8").Value = "New Transmittal from Automation"
$ws.Range("D8").Value = "UnTick"
$ws.Range("E8").Value = "Change Note"
$ws.Range("F8").Value = "Issued for Approval"
$ws.Range("L8").Value = "Message for New transmittal"
$ws.Range("M8").Value = "Rejected"
